$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'63.093.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Formula = "'  -1.23%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Formula = "'2.575.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Formula = "'  -2.81%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Formula = "'  +0.02%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Formula = "'587.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Formula = "'  -3.59%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Formula = "'149.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Formula = "'  +0.99%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Formula = "'  +0.01%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Formula = "'  -0.79%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Formula = "'0.109"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Formula = "'  +0.05%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Formula = "'5.63"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Formula = "'  +1.18%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Formula = "'0.382"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Formula = "'  +0.14%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Formula = "'  -0.69%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Formula = "'27.32"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Formula = "'  -0.77%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Formula = "'3.036.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Formula = "'  -2.72%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Formula = "'62.900.43"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Formula = "'  -1.25%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Formula = "'  +5.12%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Formula = "'2.589.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Formula = "'  -1.83%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Formula = "'12.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Formula = "'  +4.57%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Formula = "'4.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Formula = "'  +2.92%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Formula = "'344.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Formula = "'  -0.55%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Formula = "'6.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Formula = "'  -0.67%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Formula = "'  -0.13%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Formula = "'67.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Formula = "'  +1.23%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Formula = "'1.68"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Formula = "'  +1.90%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Formula = "'9.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Formula = "'  +0.70%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Formula = "'1.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Formula = "'  -2.04%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Formula = "'551.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Formula = "'  -2.24%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Formula = "'1.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Formula = "'  +1.14%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Formula = "'  -3.11%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Formula = "'  -0.97%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Formula = "'2.02"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Formula = "'  -1.28%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Formula = "'0.0₃0841"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Formula = "'  -1.73%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Formula = "'1.73"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Formula = "'  -1.68%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Formula = "'  -2.23%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Formula = "'167.71"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Formula = "'  -0.53%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Formula = "'0.409"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Formula = "'  +1.03%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Formula = "'0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Formula = "'  -0.02%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Formula = "'19.43"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Formula = "'  +1.31%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Formula = "'1.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Formula = "'  -1.36%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Formula = "'0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Formula = "'  +0.03%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Formula = "'165.97"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Formula = "'  +0.06%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Formula = "'39.51"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Formula = "'  -1.37%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Formula = "'3.92"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Formula = "'  +3.27%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Formula = "'  +2.01%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Formula = "'22.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Formula = "'  +0.74%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Formula = "'  -0.46%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Formula = "'0.0250"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Formula = "'  +2.12%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Formula = "'2.02"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Formula = "'  -0.39%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Formula = "'0.0960"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Formula = "'  -0.02%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Formula = "'18.90"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Formula = "'  +0.20%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Formula = "'0.0₆0232"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Formula = "'  +17.00%  "
$ws.Range("E51").Style = "Normal"
